# Adapt column header formatting to respective input file names (#7)
#
# - Rename the "_old" / "_new" header-name suffixes to the concrete
#   format-version suffixes "_FV2404" / "_FV2410".
# - Wrap the data range in a proper Excel Table ("Table1") with an
#   AutoFilter so the new headers are backed by a ListObject.
# - Freeze the header row (row 1) so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Base column names (without suffix), in left-to-right order as they
# appear for the "old" (FV2404) block (columns A-J) and, identically,
# for the "new" (FV2410) block (columns L-U). Column K just holds the
# literal "diff" label and is left untouched.
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

$oldCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$newCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Range($oldCols[$i] + "1").Value = $baseNames[$i] + "_FV2404"
    $ws.Range($newCols[$i] + "1").Value = $baseNames[$i] + "_FV2410"
}

# Turn the A1:U81 range into a real Excel Table ("Table1") with an
# AutoFilter, using the header row we just wrote for the column names.
$tableRange = $ws.Range("A1:U81")
$table = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$table.Name = "Table1"

# Freeze the header row (row 1) so it stays visible while scrolling.
$ws.Rows.Item(2).Select() | Out-Null
($excel.ActiveWindow.FreezePanes = $true) | Out-Null
